# Update with Correct Forecast output
#
# "Forecast Comparison" sheet: insert a new "Week_Start_Date" column
# (B), renumber the Week labels (W01 -> W1, etc.), refresh the
# MyForecast values, and store is_holiday_week as a boolean.
#
# "Summary" sheet: refresh the derived forecast totals/max that moved
# as a result of the corrected forecast numbers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: make room for the new column -------------
# Inserting before column B shifts ASIN..is_holiday_week one column to
# the right (C..J) and keeps Week in column A.
$ws1.Columns.Item(2).EntireColumn.Insert()

$ws1.Range("B1").Value = "Week_Start_Date"

$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")

$weekStartDates = @(
  "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
  "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
  "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
  "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)

$myForecast = @(675,698,684,680,523,558,630,730,602,537,558,703,504,576,514,616)

# Keep the new date column as literal text (e.g. "2025-01-05"), not an
# auto-converted date serial number.
$ws1.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
  $row = $i + 2
  $ws1.Cells.Item($row, 1).Value  = $weeks[$i]          # A: Week
  $ws1.Cells.Item($row, 2).Value  = $weekStartDates[$i] # B: Week_Start_Date
  $ws1.Cells.Item($row, 4).Value  = $myForecast[$i]     # D: MyForecast
  $ws1.Cells.Item($row, 10).Value = $false              # J: is_holiday_week (boolean)
}

# --- Summary: refresh totals derived from the corrected forecast ----
$ws2.Range("B9:B12").NumberFormat = "@"
$ws2.Range("B9").Value  = "9786"  # Total Forecast (16 Weeks)
$ws2.Range("B10").Value = "5176"  # Total Forecast (8 Weeks)
$ws2.Range("B11").Value = "2736"  # Total Forecast (4 Weeks)
$ws2.Range("B12").Value = "730"   # Max Forecast
